# Natmi following Dr Hou advice
# Rebuilds the LR-pair result rows on Sheet1 (columns A:T) with the updated
# ligand/receptor statistics, and appends three new rows (8-10) so every
# combination of Sending cluster x Target cluster (ECs / FAPs / sCs) for the
# Col1a2 -> Itga2 pair is represented.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ r = 2;  vals = @("ECs","Col1a2","Itga2","ECs", 3,1,5.256564333333333,15.769693,0.003747859920520347,0.003747859920520347,2,0.6666666666666666,1.242496666666667,3.72749,0.3272238221337332,0.3272238221337332,6.531263662285555,58.78137296057,0.001226389048014498,0.001226389048014498) },
    @{ r = 3;  vals = @("ECs","Col1a2","Itga2","FAPs",3,1,5.256564333333333,15.769693,0.003747859920520347,0.003747859920520347,3,1,1.744414,5.233242000000001,0.459408730644692,0.459408730644692,9.169624414967334,82.52661973470602,0.001721799568720369,0.001721799568720369) },
    @{ r = 4;  vals = @("ECs","Col1a2","Itga2","sCs", 3,1,5.256564333333333,15.769693,0.003747859920520347,0.003747859920520347,3,1,0.8101743333333333,2.430523,0.2133674472215748,0.2133674472215748,4.258733504382112,38.328601539439,0.0007996713037854807,0.0007996713037854804) },
    @{ r = 5;  vals = @("FAPs","Col1a2","Itga2","ECs", 3,1,1312.703450666667,3938.110352,0.93593996730609,0.9359399673060897,2,0.6666666666666666,1.242496666666667,3.72749,0.3272238221337332,0.3272238221337332,1631.029661775164,14679.26695597648,0.3062618533896201,0.30626185338962) },
    @{ r = 6;  vals = @("FAPs","Col1a2","Itga2","FAPs",3,1,1312.703450666667,3938.110352,0.93593996730609,0.9359399673060897,3,1,1.744414,5.233242000000001,0.459408730644692,0.459408730644692,2289.898277191243,20609.08449472118,0.4299789923397253,0.4299789923397252) },
    @{ r = 7;  vals = @("FAPs","Col1a2","Itga2","sCs", 3,1,1312.703450666667,3938.110352,0.93593996730609,0.9359399673060897,3,1,0.8101743333333333,2.430523,0.2133674472215748,0.2133674472215748,1063.518643008233,9571.667787074095,0.1996991215767446,0.1996991215767445) },
    @{ r = 8;  vals = @("sCs","Col1a2","Itga2","ECs", 3,1,84.59089266666666,253.772678,0.06031217277338979,0.06031217277338978,2,0.6666666666666666,1.242496666666667,3.72749,0.3272238221337332,0.3272238221337332,105.1039021686911,945.9351195182199,0.01973557969609869,0.01973557969609869) },
    @{ r = 9;  vals = @("sCs","Col1a2","Itga2","FAPs",3,1,84.59089266666666,253.772678,0.06031217277338979,0.06031217277338978,3,1,1.744414,5.233242000000001,0.459408730644692,0.459408730644692,147.5615374402307,1328.053836962076,0.02770793873624636,0.02770793873624635) },
    @{ r = 10; vals = @("sCs","Col1a2","Itga2","sCs", 3,1,84.59089266666666,253.772678,0.06031217277338979,0.06031217277338978,3,1,0.8101743333333333,2.430523,0.2133674472215748,0.2133674472215748,68.53337007228822,616.800330650594,0.01286865434104475,0.01286865434104474) }
)

foreach ($row in $rows) {
    $rnum = $row.r
    $vals = $row.vals
    for ($c = 0; $c -lt $vals.Length; $c++) {
        $ws.Cells.Item($rnum, $c + 1).Value = $vals[$c]
    }
}
